# Update result from 001 - 004
# Append 15 new log rows (rows 12-26) to the "logs" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features  = "10 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii"
$modelType = "Neural-Network"

$model1000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$model2000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$model3000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000"

$rows = @(
    @("20160426_092010", $model1000, 0.920792079207921, 0.56),
    @("20160426_093448", $model1000, 0.920792079207921, 0.45),
    @("20160426_095000", $model1000, 0.924092409240924, 0.47),
    @("20160426_100456", $model1000, 0.920792079207921, 0.47),
    @("20160426_101945", $model1000, 0.917491749174917, 0.45),
    @("20160426_110626", $model2000, 0.920792079207921, 0.45),
    @("20160426_113404", $model2000, 0.917491749174917, 0.45),
    @("20160426_120242", $model2000, 0.920792079207921, 0.45),
    @("20160426_123026", $model2000, 0.914191419141914, 0.43),
    @("20160426_125821", $model2000, 0.920792079207921, 0.43),
    @("20160426_134421", $model3000, 0.910891089108911, 0.44),
    @("20160426_142730", $model3000, 0.910891089108911, 0.45),
    @("20160426_151016", $model3000, 0.917491749174917, 0.46),
    @("20160426_155310", $model3000, 0.917491749174917, 0.46),
    @("20160426_163451", $model3000, 0.917491749174917, 0.45)
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $time = $rows[$i][0]
    $model = $rows[$i][1]
    $classifyAcc = $rows[$i][2]
    $segmentAcc = $rows[$i][3]

    $ws.Cells.Item($r, 1).Value = $time
    $ws.Cells.Item($r, 2).Value = $features
    $ws.Cells.Item($r, 3).Value = $features
    $ws.Cells.Item($r, 4).Value = $features
    $ws.Cells.Item($r, 5).Value = $modelType
    $ws.Cells.Item($r, 6).Value = $model
    $ws.Cells.Item($r, 7).Value = $modelType
    $ws.Cells.Item($r, 8).Value = $model
    $ws.Cells.Item($r, 9).Value = $modelType
    $ws.Cells.Item($r, 10).Value = $model
    $ws.Cells.Item($r, 11).Value = $classifyAcc
    $ws.Cells.Item($r, 12).Value = $segmentAcc
}
